# Added child deletion flow
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header "Aspiration" -> "Aspirations"
$ws.Range("E1").Value = "Aspirations"

# Move the active selection to E2
$ws.Range("E2").Select()
